$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (new values, previously in row 3)
$ws.Range("D2").Value = 44174
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 19000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 19500
$ws.Range("R2").Value = "Región Metropolitana"
$ws.Range("S2").Value = 1083

# Row 3 (new values, previously in row 4)
$ws.Range("D3").Value = 44169
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 21000
$ws.Range("R3").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S3").Value = 1167

# Row 4 (new values, previously in row 2)
$ws.Range("D4").Value = 44160
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 24000
$ws.Range("O4").Value = 25000
$ws.Range("P4").Value = 24500
$ws.Range("R4").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S4").Value = 1361
